# Apply updated dSF (column F) values for the listed rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 7
    12 = -5
    14 = -3
    15 = -3
    17 = 1
    19 = -2
    20 = -7
    23 = 6
    26 = -3
    31 = -6
    32 = -4
    33 = 3
    37 = 11
    38 = -4
    40 = 0
    42 = -3
    45 = -5
    46 = -3
    47 = -2
    50 = -1
    52 = 0
    54 = 5
    60 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
